$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# A reusable scratch cell holding a generic "Text" number format ("@"),
# used to flip a target cell to text-mode without Excel auto-detecting
# date-like strings as real dates.
$textFmtCell = $ws.Range("ZZ2")
$textFmtCell.NumberFormat = "@"

# A reusable scratch cell used to stash/restore a cell's original format.
$scratchCell = $ws.Range("ZZ1")

# Writes $text into $addr as a literal text value (not auto-converted to a
# date/number) while preserving the cell's original number format/style.
function Set-LiteralText($addr, $text) {
    $cell = $ws.Range($addr)

    # remember current formatting
    $cell.Copy() | Out-Null
    $scratchCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    # switch to text format so the value isn't re-interpreted, then write it
    $textFmtCell.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null           # xlPasteFormats
    $cell.Value2 = $text

    # restore the original formatting (value is untouched by a formats-only paste)
    $scratchCell.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null           # xlPasteFormats
    $scratchCell.ClearContents()
}

# Marks $addr as "Complete" using the same look as the other Complete cells
# (bold, green font) by copying the format from B2, an existing Complete cell.
function Set-CompleteStatus($addr) {
    $cell = $ws.Range($addr)
    $ws.Range("B2").Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null           # xlPasteFormats
    $cell.Value2 = "Complete"
}

# ---------------------------------------------------------------------------
# Rows 29-31 - Members Read/Create/Update -> Complete, with updated dates
# ---------------------------------------------------------------------------
Set-CompleteStatus "B29"
Set-LiteralText "D29" "05/21/2024"

Set-CompleteStatus "B30"
Set-LiteralText "D30" "05/22/2024"

Set-CompleteStatus "B31"
Set-LiteralText "D31" "05/24/2024"

# ---------------------------------------------------------------------------
# Row 13 - Task 12: Login Backend -> Complete, last-updated 05/27/2024
# ---------------------------------------------------------------------------
Set-CompleteStatus "B13"
Set-LiteralText "D13" "05/27/2024"

# ---------------------------------------------------------------------------
# New tasks 41 & 42
# ---------------------------------------------------------------------------
$ws.Range("A41").Value2 = "Task 41: Inculcate error texts in forms, add project will be ur basis since it works"
$ws.Range("A42").Value2 = "Task 42: Fix bug not closing modal upon form completion"

# ---------------------------------------------------------------------------
# Cleanup scratch cells, fix selection to match the saved view
# ---------------------------------------------------------------------------
$scratchCell.Clear()
$textFmtCell.Clear()

$ws.Range("C41").Select()
